$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed symbol snapshot.
# A leading apostrophe keeps each value as literal text (preserving trailing
# zeros and the "%" suffix) instead of Excel auto-converting it to a number.

$ws.Range("D2").Value = "'316.93"
$ws.Range("E2").Value = "'4.08%"
$ws.Range("D3").Value = "'36.10"
$ws.Range("E3").Value = "'0.34%"
$ws.Range("D4").Value = "'5.146"
$ws.Range("E4").Value = "'1.10%"
$ws.Range("D5").Value = "'0.08252"
$ws.Range("E5").Value = "'5.04%"
$ws.Range("D6").Value = "'2.151"
$ws.Range("E6").Value = "'1.69%"
$ws.Range("D7").Value = "'8.025"
$ws.Range("E7").Value = "'0.95%"
$ws.Range("D8").Value = "'0.9275"
$ws.Range("E8").Value = "'0.69%"
$ws.Range("D9").Value = "'0.1024"
$ws.Range("E9").Value = "'5.13%"
$ws.Range("D10").Value = "'0.1885"
$ws.Range("E10").Value = "'1.85%"
$ws.Range("D11").Value = "'0.09336"
$ws.Range("D12").Value = "'0.03614"
$ws.Range("E12").Value = "'1.81%"
$ws.Range("D13").Value = "'0.09918"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("D14").Value = "'0.001436"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.005719"
$ws.Range("E15").Value = "'0.32%"
$ws.Range("D16").Value = "'3.456"
$ws.Range("E16").Value = "'-0.21%"
$ws.Range("D17").Value = "'4.141"
$ws.Range("E17").Value = "'0.89%"
$ws.Range("D18").Value = "'2.798"
$ws.Range("E18").Value = "'13.49%"
$ws.Range("D19").Value = "'0.3375"
$ws.Range("E19").Value = "'-1.39%"
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'1.42%"
$ws.Range("D21").Value = "'5.189"
$ws.Range("E21").Value = "'-1.66%"
$ws.Range("D22").Value = "'0.2192"
$ws.Range("E22").Value = "'-0.46%"
$ws.Range("D23").Value = "'0.04594"
$ws.Range("E23").Value = "'0.90%"
$ws.Range("E24").Value = "'0.84%"
$ws.Range("D25").Value = "'0.004736"
$ws.Range("E25").Value = "'-6.82%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-21.93%"
$ws.Range("D27").Value = "'0.0004503"
$ws.Range("E27").Value = "'-5.35%"
$ws.Range("D39").Value = "'0.02002"
$ws.Range("E39").Value = "'8.75%"
$ws.Range("D40").Value = "'0.04947"
$ws.Range("E40").Value = "'4.73%"
$ws.Range("D41").Value = "'0.007794"
$ws.Range("E41").Value = "'3.49%"
$ws.Range("E42").Value = "'0.07%"
$ws.Range("D43").Value = "'0.007827"
$ws.Range("E43").Value = "'1.13%"
$ws.Range("D44").Value = "'0.002142"
$ws.Range("E44").Value = "'-3.42%"
$ws.Range("D45").Value = "'0.01173"
$ws.Range("E45").Value = "'6.51%"
$ws.Range("D46").Value = "'0.00006469"
$ws.Range("E46").Value = "'2.11%"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'33.67"
$ws.Range("E48").Value = "'-18.48%"
$ws.Range("D49").Value = "'0.001902"
$ws.Range("E49").Value = "'-5.09%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.20%"
